$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item(4)
$shp.TextFrame.TextRange.Text = "C"
